$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.411649666666666
$ws.Range("H2").Value = 16.234949
$ws.Range("I2").Value = 0.1787865280277313
$ws.Range("J2").Value = 0.1787865280277313
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.348457333333333
$ws.Range("N2").Value = 16.045372
$ws.Range("O2").Value = 0.4996811083875221
$ws.Range("P2").Value = 0.499681108387522
$ws.Range("Q2").Value = 28.94397734511422
$ws.Range("R2").Value = 260.495796106028
$ws.Range("S2").Value = 0.08933625048965355
$ws.Range("T2").Value = 0.08933625048965353

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.411649666666666
$ws.Range("H3").Value = 16.234949
$ws.Range("I3").Value = 0.1787865280277313
$ws.Range("J3").Value = 0.1787865280277313
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.355284000000001
$ws.Range("N3").Value = 16.065852
$ws.Range("O3").Value = 0.500318891612478
$ws.Range("P3").Value = 0.5003188916124779
$ws.Range("Q3").Value = 28.98092087350534
$ws.Range("R3").Value = 260.8282878615481
$ws.Range("S3").Value = 0.08945027753807776
$ws.Range("T3").Value = 0.08945027753807773

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.999757
$ws.Range("H4").Value = 5.999271
$ws.Range("I4").Value = 0.06606665858866914
$ws.Range("J4").Value = 0.06606665858866914
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.348457333333333
$ws.Range("N4").Value = 16.045372
$ws.Range("O4").Value = 0.4996811083875221
$ws.Range("P4").Value = 0.499681108387522
$ws.Range("Q4").Value = 10.69561499153467
$ws.Range("R4").Value = 96.260534923812
$ws.Range("S4").Value = 0.0330122611910462
$ws.Range("T4").Value = 0.03301226119104619

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.999757
$ws.Range("H5").Value = 5.999271
$ws.Range("I5").Value = 0.06606665858866914
$ws.Range("J5").Value = 0.06606665858866914
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.355284000000001
$ws.Range("N5").Value = 16.065852
$ws.Range("O5").Value = 0.500318891612478
$ws.Range("P5").Value = 0.5003188916124779
$ws.Range("Q5").Value = 10.709266665988
$ws.Range("R5").Value = 96.38339999389203
$ws.Range("S5").Value = 0.03305439739762295
$ws.Range("T5").Value = 0.03305439739762294

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.42498166666667
$ws.Range("H6").Value = 37.274945
$ws.Range("I6").Value = 0.410488385209873
$ws.Range("J6").Value = 0.410488385209873
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.348457333333333
$ws.Range("N6").Value = 16.045372
$ws.Range("O6").Value = 0.4996811083875221
$ws.Range("P6").Value = 0.499681108387522
$ws.Range("Q6").Value = 66.45448431161554
$ws.Range("R6").Value = 598.09035880454
$ws.Range("S6").Value = 0.2051132913018734
$ws.Range("T6").Value = 0.2051132913018734

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.42498166666667
$ws.Range("H7").Value = 37.274945
$ws.Range("I7").Value = 0.410488385209873
$ws.Range("J7").Value = 0.410488385209873
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.355284000000001
$ws.Range("N7").Value = 16.065852
$ws.Range("O7").Value = 0.500318891612478
$ws.Range("P7").Value = 0.5003188916124779
$ws.Range("Q7").Value = 66.53930551979334
$ws.Range("R7").Value = 598.85374967814
$ws.Range("S7").Value = 0.2053750939079995
$ws.Range("T7").Value = 0.2053750939079995

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.235186
$ws.Range("H8").Value = 27.705558
$ws.Range("I8").Value = 0.305106010612718
$ws.Range("J8").Value = 0.305106010612718
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.348457333333333
$ws.Range("N8").Value = 16.045372
$ws.Range("O8").Value = 0.4996811083875221
$ws.Range("P8").Value = 0.499681108387522
$ws.Range("Q8").Value = 49.39399828639733
$ws.Range("R8").Value = 444.545984577576
$ws.Range("S8").Value = 0.152455709558658
$ws.Range("T8").Value = 0.152455709558658

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.235186
$ws.Range("H9").Value = 27.705558
$ws.Range("I9").Value = 0.305106010612718
$ws.Range("J9").Value = 0.305106010612718
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.355284000000001
$ws.Range("N9").Value = 16.065852
$ws.Range("O9").Value = 0.500318891612478
$ws.Range("P9").Value = 0.5003188916124779
$ws.Range("Q9").Value = 49.45704382282401
$ws.Range("R9").Value = 445.1133944054161
$ws.Range("S9").Value = 0.15265030105406
$ws.Range("T9").Value = 0.15265030105406

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.197203333333333
$ws.Range("H10").Value = 3.59161
$ws.Range("I10").Value = 0.03955241756100866
$ws.Range("J10").Value = 0.03955241756100866
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.348457333333333
$ws.Range("N10").Value = 16.045372
$ws.Range("O10").Value = 0.4996811083875221
$ws.Range("P10").Value = 0.499681108387522
$ws.Range("Q10").Value = 6.403190947657778
$ws.Range("R10").Value = 57.62871852892
$ws.Range("S10").Value = 0.0197635958462909
$ws.Range("T10").Value = 0.0197635958462909

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.197203333333333
$ws.Range("H11").Value = 3.59161
$ws.Range("I11").Value = 0.03955241756100866
$ws.Range("J11").Value = 0.03955241756100866
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.355284000000001
$ws.Range("N11").Value = 16.065852
$ws.Range("O11").Value = 0.500318891612478
$ws.Range("P11").Value = 0.5003188916124779
$ws.Range("Q11").Value = 6.411363855746668
$ws.Range("R11").Value = 57.70227470172001
$ws.Range("S11").Value = 0.01978882171471777
$ws.Range("T11").Value = 0.01978882171471776
